$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 468) holds the "Förändrad" (last changed) date.
# It is being bumped from serial date 45179 (2023-09-10) to 45180 (2023-09-11)
# for every data row in the sheet.
$ws.Range("C2:C468").Value = 45180
